$d = $word.ActiveDocument

# 1. Simplify the getBalance() signature line: merges the three runs
#    ("+ getBalance()", ": ", "Double") - which also straddle the old
#    "_GoBack" bookmark - into a single run reading
#    "+ getBalance() -> Double". The Find/Replace removes the old
#    bookmark along with the replaced text, matching the diff.
$d.Content.Find.Execute("+ getBalance(): Double", $false, $false, $false, $false, $false, $true, 1, $false, "+ getBalance() -> Double", 2)

# 2. Word re-drops the "_GoBack" bookmark at the location of the most
#    recent edit - here, right after ", but valid" (before " percentage)")
#    in the "Setting discount on the card" bullet.
$rng = $d.Content
$rng.Find.Execute(", but valid")
$target = $rng.Duplicate
$target.Collapse(0)
$d.Bookmarks.Add("_GoBack", $target)
